# Fix the 保險 (insurance) and 債務 (debt) sheets: the header row had been
# accidentally populated with row-2 data instead of real column names, and
# several metadata columns (property_category/category/date/legislator_name/
# legislator_id/source_file/index, plus company/species/debtor labels) were
# missing entirely. Also turn the numeric amount columns back into real
# numbers instead of text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 保險 (insurance) -> sheet3.xml
# ---------------------------------------------------------------------
$ins = $wb.Worksheets.Item("保險")

# Header row (row 1)
$ins.Cells.Item(1,2).Value = "company"
$ins.Cells.Item(1,3).Value = "name"
$ins.Cells.Item(1,4).Value = "owner"
$ins.Cells.Item(1,5).Value = "property_category"
$ins.Cells.Item(1,6).Value = "category"
$ins.Cells.Item(1,7).Value = "date"
$ins.Cells.Item(1,8).Value = "legislator_name"
$ins.Cells.Item(1,9).Value = "legislator_id"
$ins.Cells.Item(1,10).Value = "source_file"
$ins.Cells.Item(1,11).Value = "index"

# Column G (date) holds the ISO-looking literal "2013-12-30" - force text
# format first so Excel doesn't auto-convert it into a date serial number.
$ins.Range("G2:G5").NumberFormat = "@"

# Row 2 (index 77)
$ins.Cells.Item(2,2).Value = "國泰人壽"
$ins.Cells.Item(2,3).Value = "鑫添鑫終身壽險"
$ins.Cells.Item(2,4).Value = "何欣純"
$ins.Cells.Item(2,5).Value = "insurance"
$ins.Cells.Item(2,6).Value = "normal"
$ins.Cells.Item(2,7).Value = "2013-12-30"
$ins.Cells.Item(2,8).Value = "何欣純"
$ins.Cells.Item(2,9).Value = 1733
$ins.Cells.Item(2,10).Value = "tmp8e3c1"
$ins.Cells.Item(2,11).Value = 77

# Row 3 (index 78)
$ins.Cells.Item(3,2).Value = "富邦人壽"
$ins.Cells.Item(3,3).Value = "金豐樂養老保險"
$ins.Cells.Item(3,4).Value = "何欣純"
$ins.Cells.Item(3,5).Value = "insurance"
$ins.Cells.Item(3,6).Value = "normal"
$ins.Cells.Item(3,7).Value = "2013-12-30"
$ins.Cells.Item(3,8).Value = "何欣純"
$ins.Cells.Item(3,9).Value = 1733
$ins.Cells.Item(3,10).Value = "tmp8e3c1"
$ins.Cells.Item(3,11).Value = 78

# Row 4 (index 79)
$ins.Cells.Item(4,2).Value = "國寶人壽"
$ins.Cells.Item(4,3).Value = "得意年年终身壽險"
$ins.Cells.Item(4,4).Value = "何欣純"
$ins.Cells.Item(4,5).Value = "insurance"
$ins.Cells.Item(4,6).Value = "normal"
$ins.Cells.Item(4,7).Value = "2013-12-30"
$ins.Cells.Item(4,8).Value = "何欣純"
$ins.Cells.Item(4,9).Value = 1733
$ins.Cells.Item(4,10).Value = "tmp8e3c1"
$ins.Cells.Item(4,11).Value = 79

# Row 5 (index 80)
$ins.Cells.Item(5,2).Value = "國泰人壽"
$ins.Cells.Item(5,3).Value = "富貴保本三福終身壽險"
$ins.Cells.Item(5,4).Value = "謝俊雄"
$ins.Cells.Item(5,5).Value = "insurance"
$ins.Cells.Item(5,6).Value = "normal"
$ins.Cells.Item(5,7).Value = "2013-12-30"
$ins.Cells.Item(5,8).Value = "何欣純"
$ins.Cells.Item(5,9).Value = 1733
$ins.Cells.Item(5,10).Value = "tmp8e3c1"
$ins.Cells.Item(5,11).Value = 80

# ---------------------------------------------------------------------
# Sheet 債務 (debt) -> sheet4.xml
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")

# Header row (row 1)
$debt.Cells.Item(1,2).Value = "species"
$debt.Cells.Item(1,3).Value = "debtor"
$debt.Cells.Item(1,4).Value = "owner"
$debt.Cells.Item(1,5).Value = "total"
$debt.Cells.Item(1,6).Value = "register_date"
$debt.Cells.Item(1,7).Value = "register_reason"
$debt.Cells.Item(1,8).Value = "property_category"
$debt.Cells.Item(1,9).Value = "category"
$debt.Cells.Item(1,10).Value = "date"
$debt.Cells.Item(1,11).Value = "legislator_name"
$debt.Cells.Item(1,12).Value = "legislator_id"
$debt.Cells.Item(1,13).Value = "source_file"
$debt.Cells.Item(1,14).Value = "index"

# Column J (date) holds the ISO-looking literal "2013-12-30" - force text
# format first so Excel doesn't auto-convert it into a date serial number.
$debt.Range("J2:J3").NumberFormat = "@"

# Row 2 (index 90)
$debt.Cells.Item(2,2).Value = "房屋貸款"
$debt.Cells.Item(2,3).Value = "謝俊雄"
$debt.Cells.Item(2,4).Value = "霧峰鄉農會臺中市霧峰區四德路10號"
$debt.Cells.Item(2,5).Value = 5215377
$debt.Cells.Item(2,6).Value = "93年09月14日"
$debt.Cells.Item(2,7).Value = "房貸"
$debt.Cells.Item(2,8).Value = "debt"
$debt.Cells.Item(2,9).Value = "normal"
$debt.Cells.Item(2,10).Value = "2013-12-30"
$debt.Cells.Item(2,11).Value = "何欣純"
$debt.Cells.Item(2,12).Value = 1733
$debt.Cells.Item(2,13).Value = "tmp8e3c1"
$debt.Cells.Item(2,14).Value = 90

# Row 3 (index 91)
$debt.Cells.Item(3,2).Value = "房屋貸款"
$debt.Cells.Item(3,3).Value = "謝俊雄"
$debt.Cells.Item(3,4).Value = "台中商業銀行清水分行臺中市清水區中山路104號"
$debt.Cells.Item(3,5).Value = 1500000
$debt.Cells.Item(3,6).Value = "102年09月26日"
$debt.Cells.Item(3,7).Value = "房貸"
$debt.Cells.Item(3,8).Value = "debt"
$debt.Cells.Item(3,9).Value = "normal"
$debt.Cells.Item(3,10).Value = "2013-12-30"
$debt.Cells.Item(3,11).Value = "何欣純"
$debt.Cells.Item(3,12).Value = 1733
$debt.Cells.Item(3,13).Value = "tmp8e3c1"
$debt.Cells.Item(3,14).Value = 91
